$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: apply the word-level text corrections (typo fixes / rewordings).
# wdFindContinue = 1, wdReplaceAll = 2
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("séannce", $false, $false, $false, $false, $false, $true, 1, $false, "séance", 2)
$d.Content.Find.Execute("Aujourd’hu nous", $false, $false, $false, $false, $false, $true, 1, $false, "Aujourd’hui nous", 2)
$d.Content.Find.Execute("empecher", $false, $false, $false, $false, $false, $true, 1, $false, "empêcher", 2)
$d.Content.Find.Execute("rems à", $false, $false, $false, $false, $false, $true, 1, $false, "remis à", 2)
$d.Content.Find.Execute("normalemnt", $false, $false, $false, $false, $false, $true, 1, $false, "normalement", 2)
$d.Content.Find.Execute("explosions sont dû", $false, $false, $false, $false, $false, $true, 1, $false, "explosions étaient dûe", 2)
$d.Content.Find.Execute("suivante explosait", $false, $false, $false, $false, $false, $true, 1, $false, "suivantes explosaient", 2)
$d.Content.Find.Execute("soudure ne", $false, $false, $false, $false, $false, $true, 1, $false, "soudures ne", 2)
$d.Content.Find.Execute("des fls se", $false, $false, $false, $false, $false, $true, 1, $false, "des fils se", 2)
$d.Content.Find.Execute("à reflechir à", $false, $false, $false, $false, $false, $true, 1, $false, "à réfléchir à", 2)
$d.Content.Find.Execute("geométrique.", $false, $false, $false, $false, $false, $true, 1, $false, "geométriques.", 2)
$d.Content.Find.Execute("attandant", $false, $false, $false, $false, $false, $true, 1, $false, "attendant", 2)

# ---------------------------------------------------------------------------
# Step 1b: the original document had a run break (with a _GoBack bookmark in
# between) right after "...suivante explosait. Fin". In the target layout
# that break is gone and the two runs are merged into one continuous run, so
# remove the bookmark and re-save identical text across the old boundary to
# force Word to coalesce the two runs back into a single run.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$d.Content.Find.Execute("suivantes explosaient. Finalement", $false, $false, $false, $false, $false, $true, 1, $false, "suivantes explosaient. Finalement", 2)

# ---------------------------------------------------------------------------
# Step 2: re-create the run boundaries that show up in the target XML. These
# boundaries are introduced by toggling (and un-toggling) a character
# property on a growing sub-range of each paragraph; because the property
# ends up back at its original value the rendered formatting is unchanged,
# but Word is forced to break the paragraph's text into multiple <w:r> runs
# at each of those offsets (matching the shape of the diff).
# ---------------------------------------------------------------------------
function Split-Run($range_start, $offsets) {
    foreach ($o in $offsets) {
        $r = $d.Range($range_start, $range_start + $o)
        $r.Font.Bold = $true
        $r.Font.Bold = $false
    }
}

$p1Start = $d.Paragraphs(1).Range.Start
Split-Run $p1Start @(11)

$p2Start = $d.Paragraphs(2).Range.Start
Split-Run $p2Start @(10, 11, 42, 74, 273, 274, 294, 295, 299, 334, 337, 338, 419, 420, 429, 431, 563, 564, 593, 594)

$p3Start = $d.Paragraphs(3).Range.Start
Split-Run $p3Start @(29, 76, 77, 154, 336, 376, 408)

# ---------------------------------------------------------------------------
# Step 3: re-create the _GoBack bookmark at its new position (right after
# "recuper" in paragraph 3; it used to sit right after "Fin" in paragraph 2,
# see step 1b above where the old one was removed).
# ---------------------------------------------------------------------------
$p3Start = $d.Paragraphs(3).Range.Start
$newBmRange = $d.Range($p3Start + 336, $p3Start + 336)
$d.Bookmarks.Add("_GoBack", $newBmRange)

Write-Host "Edit complete"
